# Azerbaijan Premier League workbook update (19-04-2024 23:27)
# A new match result (PFK Turan Tovuz vs FK Sumqayit) is inserted as row 151,
# pushing the previously-last row (Araz FK vs Sabail FC, id 149) down to row
# 152 where it is re-numbered id 150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 151; this shifts the former row 151
# (and everything that was empty below it) down by one, so the old last
# match now lives at row 152.
$ws.Rows.Item(151).Insert()

# The new row 151 should look like the other data rows: bold/centered/
# bordered "id" cell in column A, and a date-formatted cell in column E.
$ws.Range("A151").Font.Bold = $true
$ws.Range("A151").HorizontalAlignment = -4108
$ws.Range("A151").VerticalAlignment = -4160
$ws.Range("A151").Borders.LineStyle = 1
$ws.Range("E151").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Renumber the id of the (now shifted) old last row: 149 -> 150
$ws.Range("A152").Value = 150

# Populate the new match row (id 149) with its data
$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 7020807
$ws.Range("C151").Value = "Azerbaijan Premier League"
$ws.Range("D151").Value = "Azerbaijan Premier League"
$ws.Range("E151").Value = 45401.39583333334
$ws.Range("F151").Value = "PFK Turan Tovuz"
$ws.Range("G151").Value = "FK Sumqayit"
$ws.Range("H151").Value = 1
$ws.Range("I151").Value = 4
$ws.Range("J151").Value = "A"
$ws.Range("K151").Value = 2.25
$ws.Range("L151").Value = 3.1
$ws.Range("M151").Value = 2.9
$ws.Range("N151").Value = 2
$ws.Range("O151").Value = 3
$ws.Range("P151").Value = 3.4
$ws.Range("Q151").Value = -0.25
$ws.Range("R151").Value = 1.775
$ws.Range("S151").Value = 2.025
$ws.Range("T151").Value = 2.25
$ws.Range("U151").Value = 2.025
$ws.Range("V151").Value = 1.775
$ws.Range("W151").Value = -1
$ws.Range("X151").Value = -1
$ws.Range("Y151").Value = 2.4
$ws.Range("Z151").Value = -1
$ws.Range("AA151").Value = 1.025
$ws.Range("AB151").Value = 1.025
$ws.Range("AC151").Value = -1

Write-Output ("New used range: " + $ws.UsedRange.Address())
